$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Negative Case update on "Pelanggan" (Customer) test data: edit F2 ---
$ws.Range("F2").Value = "Edited Rotten Sugar"

# --- New test case rows: Bahan Kemas, Bahan Pembantu, Barang Setengah Jadi ---

# Row 3: stray copy of kodebarang value (GB01) in column B, carrying the same
# formatting as the header/row2 block.
$ws.Range("B3").Value = "GB01"
$ws.Range("B3").Interior.ColorIndex = -4142

# Row 4: full duplicate of the "Gula Busuk" record (positive case), same
# formatting as row 2.
$ws.Range("A4").Value = "Gula Busuk"
$ws.Range("B4").Value = "GB01"
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = "Kilogram"
$ws.Range("A4:D4").Interior.ColorIndex = -4142

# Row 5: Bahan Kemas negative-case row (Daun Teh Busuk) - default formatting.
$ws.Range("A5").Value = "Daun Teh Busuk"
$ws.Range("C5").Value = "Qw2E0#"
$ws.Range("D5").Value = "Kodi"

# Row 6: Bahan Pembantu negative-case row (Kakao Busuk) - distinct formatting.
$ws.Range("A6").Value = "Kakao Busuk"
$ws.Range("C6").Value = "ASD"
$ws.Range("D6").Value = "Kilogram"
$ws.Range("A6:D6").Interior.ColorIndex = -4142

# Row 7: Barang Setengah Jadi test case (Singkong Busuk) plus a copied
# "Gula Busuk" reference block in F7:I7 (same formatting as row 2).
$ws.Range("A7").Value = "Singkong Busuk"
$ws.Range("B7").Value = "SK01"
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = "Kilogram"

$ws.Range("F7").Value = "Gula Busuk"
$ws.Range("G7").Value = "GB01"
$ws.Range("H7").Value = 1
$ws.Range("I7").Value = "Kilogram"
$ws.Range("F7:I7").Interior.ColorIndex = -4142

# Column A is now wide enough for "Daun Teh Busuk" / "Singkong Busuk".
$ws.Columns.Item(1).ColumnWidth = 14.75

# Final selection left on the newly pasted reference block.
$ws.Range("F7:I7").Select()
